# Update gh-pages to output generated at 456a3b4
# Applies the refreshed "remaining stock" (F column) numbers and marks a
# few listings as sold out ("已售罄") in the G column, across the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 510
$ws1.Range("F4").Value  = 489
$ws1.Range("F6").Value  = 142
$ws1.Range("F7").Value  = 902
$ws1.Range("F8").Value  = 704
$ws1.Range("F9").Value  = 158
$ws1.Range("F11").Value = 75
$ws1.Range("F12").Value = 754
$ws1.Range("F13").Value = 245
$ws1.Range("F14").Value = 538
$ws1.Range("F16").Value = 1270
$ws1.Range("F17").Value = 110
$ws1.Range("F18").Value = 299
$ws1.Range("F19").Value = 1018
$ws1.Range("F20").Value = 2758
$ws1.Range("F21").Value = 1234
$ws1.Range("F23").Value = 160
$ws1.Range("F24").Value = 1227
$ws1.Range("F26").Value = 950
$ws1.Range("F27").Value = 306
$ws1.Range("F28").Value = 79
$ws1.Range("F29").Value = 1281

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("G3").Value = "已售罄"
$ws2.Range("F7").Value = 2

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 713
$ws4.Range("F3").Value  = 510
$ws4.Range("F6").Value  = 489
$ws4.Range("G7").Value  = "已售罄"
$ws4.Range("G8").Value  = "已售罄"
$ws4.Range("F13").Value = 142
$ws4.Range("F14").Value = 2
$ws4.Range("F15").Value = 902
$ws4.Range("F16").Value = 704
$ws4.Range("F17").Value = 158
$ws4.Range("F23").Value = 75
$ws4.Range("F25").Value = 754
$ws4.Range("F26").Value = 245
$ws4.Range("F27").Value = 538
$ws4.Range("F29").Value = 1270
$ws4.Range("F30").Value = 110
$ws4.Range("F31").Value = 299
$ws4.Range("F32").Value = 1018
$ws4.Range("F33").Value = 2758
$ws4.Range("F34").Value = 1234
$ws4.Range("F36").Value = 160
$ws4.Range("F37").Value = 1227
$ws4.Range("F40").Value = 950
$ws4.Range("F41").Value = 306
$ws4.Range("F42").Value = 80
$ws4.Range("F43").Value = 1281
